{"js": "// Update the date in the title paragraph and the 25 division-problem\n// answers in the table body, per the commit diff. Each text value is\n// replaced in-place (via a scoped search + insertText \"Replace\") so\n// that run/paragraph formatting (fonts, size, alignment) is preserved.\n\n// 1) Title line: date changes from 2025-08-29 Friday -> 2025-08-30 Saturday.\nconst titleResults = context.document.body.search(\"2025-08-29 Friday\", { matchCase: true });\ntitleResults.load(\"items\");\nawait context.sync();\ntitleResults.items[0].insertText(\"2025-08-30 Saturday\", Word.InsertLocation.replace);\n\n// 2) Table cell replacements. The table has 20 rows (5 \"data\" rows with\n// content, each followed by 3 blank spacer rows) and 5 columns.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// (row, col, oldText, newText) for every populated cell, in document order.\nconst cellEdits = [\n  [0, 0, \"11\u00f72=5, 1\", \"30\u00f75=6, 0\"],\n  [0, 1, \"99\u00f79=11, 0\", \"96\u00f79=10, 6\"],\n  [0, 2, \"47\u00f74=11, 3\", \"75\u00f77=10, 5\"],\n  [0, 3, \"54\u00f73=18, 0\", \"28\u00f76=4, 4\"],\n  [0, 4, \"51\u00f73=17, 0\", \"82\u00f72=41, 0\"],\n\n  [4, 0, \"15\u00f72=7, 1\", \"63\u00f79=7, 0\"],\n  [4, 1, \"83\u00f75=16, 3\", \"55\u00f75=11, 0\"],\n  [4, 2, \"19\u00f76=3, 1\", \"39\u00f78=4, 7\"],\n  [4, 3, \"24\u00f73=8, 0\", \"64\u00f74=16, 0\"],\n  [4, 4, \"13\u00f72=6, 1\", \"11\u00f75=2, 1\"],\n\n  [8, 0, \"85\u00f77=12, 1\", \"71\u00f77=10, 1\"],\n  [8, 1, \"76\u00f78=9, 4\", \"47\u00f73=15, 2\"],\n  [8, 2, \"48\u00f72=24, 0\", \"62\u00f74=15, 2\"],\n  [8, 3, \"28\u00f73=9, 1\", \"20\u00f72=10, 0\"],\n  [8, 4, \"15\u00f77=2, 1\", \"29\u00f75=5, 4\"],\n\n  [12, 0, \"56\u00f74=14, 0\", \"62\u00f74=15, 2\"],\n  [12, 1, \"30\u00f76=5, 0\", \"31\u00f77=4, 3\"],\n  [12, 2, \"40\u00f73=13, 1\", \"77\u00f77=11, 0\"],\n  [12, 3, \"50\u00f79=5, 5\", \"31\u00f75=6, 1\"],\n  [12, 4, \"23\u00f78=2, 7\", \"11\u00f76=1, 5\"],\n\n  [16, 0, \"31\u00f75=6, 1\", \"53\u00f77=7, 4\"],\n  [16, 1, \"97\u00f78=12, 1\", \"84\u00f73=28, 0\"],\n  [16, 2, \"16\u00f73=5, 1\", \"49\u00f78=6, 1\"],\n  [16, 3, \"56\u00f79=6, 2\", \"90\u00f76=15, 0\"],\n  [16, 4, \"43\u00f73=14, 1\", \"82\u00f75=16, 2\"],\n];\n\n// Scope each search to its own cell body so that cells whose new value\n// happens to equal another cell's old value (e.g. \"31\u00f75=6, 1\") cannot be\n// matched/replaced incorrectly.\nconst cellSearches = [];\nfor (const [row, col, oldText] of cellEdits) {\n  const cell = table.getCell(row, col);\n  const found = cell.body.search(oldText, { matchCase: true });\n  found.load(\"items\");\n  cellSearches.push(found);\n}\nawait context.sync();\n\nfor (let i = 0; i < cellEdits.length; i++) {\n  const [, , , newText] = cellEdits[i];\n  cellSearches[i].items[0].insertText(newText, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Update the date in the title paragraph and the 25 division-problem\n# answers in the table body, per the commit diff. Each cell's Range.Text\n# is assigned directly (by row/col position) so that run/paragraph\n# formatting (fonts, size, alignment) is preserved and there is no risk\n# of one cell's new value colliding with another cell's still-unprocessed\n# old value (several of the answers happen to reuse the same digits).\n\n$d = $word.ActiveDocument\n\n# 1) Title line: date changes from 2025-08-29 Friday -> 2025-08-30 Saturday.\n$titleRange = $d.Content\n$found = $titleRange.Find.Execute(\"2025-08-29 Friday\")\nif ($found) {\n    $titleRange.Text = \"2025-08-30 Saturday\"\n}\n\n# 2) Table cell replacements. The table has 20 rows (5 \"data\" rows with\n# content, each followed by 3 blank spacer rows) and 5 columns. Word's\n# Table.Cell(row, col) is 1-indexed.\n$table = $d.Tables.Item(1)\n\n# (row, col, newText) for every populated cell, in document order.\n$cellEdits = @(\n    @(1, 1, \"30\u00f75=6, 0\"),\n    @(1, 2, \"96\u00f79=10, 6\"),\n    @(1, 3, \"75\u00f77=10, 5\"),\n    @(1, 4, \"28\u00f76=4, 4\"),\n    @(1, 5, \"82\u00f72=41, 0\"),\n\n    @(5, 1, \"63\u00f79=7, 0\"),\n    @(5, 2, \"55\u00f75=11, 0\"),\n    @(5, 3, \"39\u00f78=4, 7\"),\n    @(5, 4, \"64\u00f74=16, 0\"),\n    @(5, 5, \"11\u00f75=2, 1\"),\n\n    @(9, 1, \"71\u00f77=10, 1\"),\n    @(9, 2, \"47\u00f73=15, 2\"),\n    @(9, 3, \"62\u00f74=15, 2\"),\n    @(9, 4, \"20\u00f72=10, 0\"),\n    @(9, 5, \"29\u00f75=5, 4\"),\n\n    @(13, 1, \"62\u00f74=15, 2\"),\n    @(13, 2, \"31\u00f77=4, 3\"),\n    @(13, 3, \"77\u00f77=11, 0\"),\n    @(13, 4, \"31\u00f75=6, 1\"),\n    @(13, 5, \"11\u00f76=1, 5\"),\n\n    @(17, 1, \"53\u00f77=7, 4\"),\n    @(17, 2, \"84\u00f73=28, 0\"),\n    @(17, 3, \"49\u00f78=6, 1\"),\n    @(17, 4, \"90\u00f76=15, 0\"),\n    @(17, 5, \"82\u00f75=16, 2\")\n)\n\nforeach ($edit in $cellEdits) {\n    $row = $edit[0]\n    $col = $edit[1]\n    $newText = $edit[2]\n    $cell = $table.Cell($row, $col)\n    $cell.Range.Text = $newText\n}\n"}
